# Add an "ImageName" column (D) of per-row image file names for the
# Kanji N3 table, matching the "Added Kanji Images for n3" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$images = @(
    "n3_feeling.jpg",
    "n3_thought.jpg",
    "n3_love.jpg",
    "n3_machine.jpg",
    "n3_pass_through.jpg",
    "n3_finish.jpg",
    "n3_emotion.webp",
    "n3_report.png",
    "n3_establish.jpg",
    "n3_manage.png"
)

for ($i = 0; $i -lt $images.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $images[$i]
}

# Page setup touched (paper size / orientation) as part of the edit.
$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait

$ws.Range("I12").Select()
